# "Generate Report for Handoff"
# Adds two new handoff rows (a .md file and a .png file) to every sheet of
# the localization-status workbook, and refreshes the existing row's
# timestamps / target-file hyperlinks to reflect a newer handoff run.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")

# Refresh existing row 2 (file renamed/re-handed-off as a .png)
$ov.Range("A2").Hyperlinks.Delete()
$ov.Hyperlinks.Add($ov.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/8b10e35ba04e41bd3e8e2f7d23b39b45f63e9b4f/e2e/27aeffc1-cf85-4d0f-bf83-a57d9d644e9d.png", "", "", "27aeffc1-cf85-4d0f-bf83-a57d9d644e9d.png")
$ov.Range("B2").Value = "Ready for handoff"
$ov.Range("C2").Value = "Ready for handoff"
$ov.Range("D2").Value = "2016-41-11 16:41:42"

# New row 3: the dependency .md file
$ov.Hyperlinks.Add($ov.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/8b10e35ba04e41bd3e8e2f7d23b39b45f63e9b4f/e2e/762bd805-4317-48f5-ac66-87a98cf5642c.md", "", "", "762bd805-4317-48f5-ac66-87a98cf5642c.md")
$ov.Range("B3").Value = "Ready for handoff"
$ov.Range("C3").Value = "Ready for handoff"
$ov.Range("D3").Value = "2016-41-11 16:41:42"

# New row 4: the new .png file
$ov.Hyperlinks.Add($ov.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/8b10e35ba04e41bd3e8e2f7d23b39b45f63e9b4f/e2e/aae16238-af30-4735-9594-f8ebebb0c6ac.png", "", "", "aae16238-af30-4735-9594-f8ebebb0c6ac.png")
$ov.Range("B4").Value = "Ready for handoff"
$ov.Range("C4").Value = "Ready for handoff"
$ov.Range("D4").Value = "2016-41-11 16:41:42"

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

# Refresh existing row 2
$zh.Range("A2").Hyperlinks.Delete()
$zh.Range("B2").Hyperlinks.Delete()
$zh.Range("D2").Hyperlinks.Delete()
$zh.Hyperlinks.Add($zh.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/8b10e35ba04e41bd3e8e2f7d23b39b45f63e9b4f/e2e/27aeffc1-cf85-4d0f-bf83-a57d9d644e9d.png", "", "", "27aeffc1-cf85-4d0f-bf83-a57d9d644e9d.png")
$zh.Hyperlinks.Add($zh.Range("B2"), "https://github.com/OpenLocalizationTest/oltest/blob/8b10e35ba04e41bd3e8e2f7d23b39b45f63e9b4f/e2e/27aeffc1-cf85-4d0f-bf83-a57d9d644e9d.png", "", "", ".png")
$zh.Range("C2").Value = "Ready for handoff"
$zh.Hyperlinks.Add($zh.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/88df66646773b65fd617eae8e87efdc98230847d/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/f71907011c90e34704f4b23e9558d9db23acaeb0.png", "", "", "f71907011c90e34704f4b23e9558d9db23acaeb0.png")
$zh.Range("E2").Value = "2016-03-11 16:41:37"
$zh.Range("H2").Value = "0001-01-01 00:00:00"
$zh.Range("I2").Value = "IsDependency"
$zh.Range("J2").Value = "e2e\762bd805-4317-48f5-ac66-87a98cf5642c.md"

# New row 3: the dependency .md file
$zh.Hyperlinks.Add($zh.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/8b10e35ba04e41bd3e8e2f7d23b39b45f63e9b4f/e2e/762bd805-4317-48f5-ac66-87a98cf5642c.md", "", "", "762bd805-4317-48f5-ac66-87a98cf5642c.md")
$zh.Hyperlinks.Add($zh.Range("B3"), "https://github.com/OpenLocalizationTest/oltest/blob/8b10e35ba04e41bd3e8e2f7d23b39b45f63e9b4f/e2e/762bd805-4317-48f5-ac66-87a98cf5642c.md", "", "", ".md")
$zh.Range("C3").Value = "Ready for handoff"
$zh.Hyperlinks.Add($zh.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/88df66646773b65fd617eae8e87efdc98230847d/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/762bd805-4317-48f5-ac66-87a98cf5642c.0257b9652a17a6000c323a3a58e794bd3d7f1aa1.zh-cn.xlf", "", "", "762bd805-4317-48f5-ac66-87a98cf5642c.0257b9652a17a6000c323a3a58e794bd3d7f1aa1.zh-cn.xlf")
$zh.Range("E3").Value = "2016-03-11 16:41:37"
$zh.Range("H3").Value = "0001-01-01 00:00:00"
$zh.Range("I3").Value = "Include"

# New row 4: the new .png file
$zh.Hyperlinks.Add($zh.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/8b10e35ba04e41bd3e8e2f7d23b39b45f63e9b4f/e2e/aae16238-af30-4735-9594-f8ebebb0c6ac.png", "", "", "aae16238-af30-4735-9594-f8ebebb0c6ac.png")
$zh.Hyperlinks.Add($zh.Range("B4"), "https://github.com/OpenLocalizationTest/oltest/blob/8b10e35ba04e41bd3e8e2f7d23b39b45f63e9b4f/e2e/aae16238-af30-4735-9594-f8ebebb0c6ac.png", "", "", ".png")
$zh.Range("C4").Value = "Ready for handoff"
$zh.Hyperlinks.Add($zh.Range("D4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/88df66646773b65fd617eae8e87efdc98230847d/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/55531ad8bc80eb4acd17c5fc42e7ac3203f6b2f0.png", "", "", "55531ad8bc80eb4acd17c5fc42e7ac3203f6b2f0.png")
$zh.Range("E4").Value = "2016-03-11 16:41:37"
$zh.Range("H4").Value = "0001-01-01 00:00:00"
$zh.Range("I4").Value = "IsDependency"
$zh.Range("J4").Value = "e2e\762bd805-4317-48f5-ac66-87a98cf5642c.md"

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

# Refresh existing row 2
$de.Range("A2").Hyperlinks.Delete()
$de.Range("B2").Hyperlinks.Delete()
$de.Range("D2").Hyperlinks.Delete()
$de.Hyperlinks.Add($de.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/8b10e35ba04e41bd3e8e2f7d23b39b45f63e9b4f/e2e/27aeffc1-cf85-4d0f-bf83-a57d9d644e9d.png", "", "", "27aeffc1-cf85-4d0f-bf83-a57d9d644e9d.png")
$de.Hyperlinks.Add($de.Range("B2"), "https://github.com/OpenLocalizationTest/oltest/blob/8b10e35ba04e41bd3e8e2f7d23b39b45f63e9b4f/e2e/27aeffc1-cf85-4d0f-bf83-a57d9d644e9d.png", "", "", ".png")
$de.Range("C2").Value = "Ready for handoff"
$de.Hyperlinks.Add($de.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/f98519993c469c1df216e50f6d9b22ede0eadac0/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/f71907011c90e34704f4b23e9558d9db23acaeb0.png", "", "", "f71907011c90e34704f4b23e9558d9db23acaeb0.png")
$de.Range("E2").Value = "2016-03-11 16:41:42"
$de.Range("H2").Value = "0001-01-01 00:00:00"
$de.Range("I2").Value = "IsDependency"
$de.Range("J2").Value = "e2e\762bd805-4317-48f5-ac66-87a98cf5642c.md"

# New row 3: the dependency .md file
$de.Hyperlinks.Add($de.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/8b10e35ba04e41bd3e8e2f7d23b39b45f63e9b4f/e2e/762bd805-4317-48f5-ac66-87a98cf5642c.md", "", "", "762bd805-4317-48f5-ac66-87a98cf5642c.md")
$de.Hyperlinks.Add($de.Range("B3"), "https://github.com/OpenLocalizationTest/oltest/blob/8b10e35ba04e41bd3e8e2f7d23b39b45f63e9b4f/e2e/762bd805-4317-48f5-ac66-87a98cf5642c.md", "", "", ".md")
$de.Range("C3").Value = "Ready for handoff"
$de.Hyperlinks.Add($de.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/f98519993c469c1df216e50f6d9b22ede0eadac0/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/762bd805-4317-48f5-ac66-87a98cf5642c.0257b9652a17a6000c323a3a58e794bd3d7f1aa1.de-de.xlf", "", "", "762bd805-4317-48f5-ac66-87a98cf5642c.0257b9652a17a6000c323a3a58e794bd3d7f1aa1.de-de.xlf")
$de.Range("E3").Value = "2016-03-11 16:41:42"
$de.Range("H3").Value = "0001-01-01 00:00:00"
$de.Range("I3").Value = "Include"

# New row 4: the new .png file
$de.Hyperlinks.Add($de.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/8b10e35ba04e41bd3e8e2f7d23b39b45f63e9b4f/e2e/aae16238-af30-4735-9594-f8ebebb0c6ac.png", "", "", "aae16238-af30-4735-9594-f8ebebb0c6ac.png")
$de.Hyperlinks.Add($de.Range("B4"), "https://github.com/OpenLocalizationTest/oltest/blob/8b10e35ba04e41bd3e8e2f7d23b39b45f63e9b4f/e2e/aae16238-af30-4735-9594-f8ebebb0c6ac.png", "", "", ".png")
$de.Range("C4").Value = "Ready for handoff"
$de.Hyperlinks.Add($de.Range("D4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/f98519993c469c1df216e50f6d9b22ede0eadac0/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/55531ad8bc80eb4acd17c5fc42e7ac3203f6b2f0.png", "", "", "55531ad8bc80eb4acd17c5fc42e7ac3203f6b2f0.png")
$de.Range("E4").Value = "2016-03-11 16:41:42"
$de.Range("H4").Value = "0001-01-01 00:00:00"
$de.Range("I4").Value = "IsDependency"
$de.Range("J4").Value = "e2e\762bd805-4317-48f5-ac66-87a98cf5642c.md"

Write-Output "Report generation complete."
